# "change tag ontology ID"
# MetaboLights - NMR sample: retag the "Tags" row from MS_1000457 (MS) to
# ARC_00000070 (ARC), and add the Comment[isObsolete] = false row that goes
# along with this template metadata refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# Row 14/15 hold the "Tags Term Accession Number" / "Tags Term Source REF"
# values for the four TAGS (columns B..E). Column C currently points at the
# MS ontology term; repoint it at the ARC ontology term.
$ws.Range("C14").Value = "http://purl.obolibrary.org/obo/ARC_00000070"
$ws.Range("C15").Value = "ARC"

# Insert a new metadata row right after the TAGS block (row 15) and before
# the AUTHORS section, carrying Comment[isObsolete] = false.
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = "Comment[isObsolete]"

# Write "false" as literal text (not a TRUE/FALSE boolean) by using Excel's
# leading-apostrophe text-prefix, then strip the resulting quote-prefix
# formatting by re-pasting the (unformatted) style from A1 over it.
$ws.Range("C16").Value = "'false"
$ws.Range("A1").Copy()
$ws.Range("C16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
